$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# 1) Update the existing paragraph's wording:
#    "Reading habits and academic success" -> "Reading habits Vs : "
$para2 = $tr.Paragraphs(2, 1)
$run1 = $para2.Runs(1, 1)
$run1.Text = "Reading habits Vs : "

# 2) Append four new bulleted paragraphs after it, matching the same
#    run formatting (sz=3600, not bold, tx1 solid fill) and bullet style
#    used elsewhere in this deck (Arial bullet character, spacing before = 0).
$bullets = "Academic performance", "Education dropout", "Level of literacy", "Education engagement"

foreach ($bulletText in $bullets) {
    $beforeLen = $tr.Length
    [void]$tr.InsertAfter("`r" + $bulletText)
    $insertedLen = $tr.Length - $beforeLen
    $newPara = $tr.Characters($beforeLen + 1, $insertedLen)

    $pf = $newPara.ParagraphFormat
    $pf.SpaceBefore = 0

    $bf = $pf.Bullet
    $bf.Font.Name = "Arial"
    $bf.Visible = -1
    $bf.Character = 8226
}
